# issue #5: property land done
#
# Cleans up stray whitespace / thousands-separators that had crept into the
# scraped text on the land/building/debt sheets, and appends the scrape
# pipeline's metadata columns (I:O -- property_category, category, date,
# legislator_name, legislator_id, source_file, index) to the "land" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1 ("土地" / land)
# ---------------------------------------------------------------------
$landSheet = $wb.Worksheets.Item(1)

# Re-key the existing data row text (strip stray internal spaces / commas
# introduced by the original scrape).
$landSheet.Range("B2").Value = "臺北市大安區龍泉段一小段02930000地號"
$landSheet.Range("D2").Value = "100000分之16216"
$landSheet.Range("F2").Value = "92年12月25日"
$landSheet.Range("G2").Value = "033貝賣"
$landSheet.Range("H2").Value = "25000000(土地建物與車位合併價）"

# New pipeline metadata columns -- copy the existing header/data formatting
# across first so the new cells pick up the same style as the rest of the
# row, then fill in the values.
$landSheet.Range("H1").Copy() | Out-Null
$landSheet.Range("I1:O1").PasteSpecial(-4122) | Out-Null
$landSheet.Range("H2").Copy() | Out-Null
$landSheet.Range("I2:O2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$landSheet.Range("I1").Value = "property_category"
$landSheet.Range("J1").Value = "category"
$landSheet.Range("K1").Value = "date"
$landSheet.Range("L1").Value = "legislator_name"
$landSheet.Range("M1").Value = "legislator_id"
$landSheet.Range("N1").Value = "source_file"
$landSheet.Range("O1").Value = "index"

$landSheet.Range("I2").Value = "land"
$landSheet.Range("J2").Value = "normal"
# Leading apostrophe forces this to stay literal text ("2012-04-30")
# instead of being auto-converted to a date serial number.
$landSheet.Range("K2").Value = "'2012-04-30"
$landSheet.Range("L2").Value = "高金素梅"
$landSheet.Range("M2").Value = 926
$landSheet.Range("N2").Value = "tmp92521"
$landSheet.Range("O2").Value = 14

# ---------------------------------------------------------------------
# Sheet 2 ("建物" / building)
# ---------------------------------------------------------------------
$buildingSheet = $wb.Worksheets.Item(2)
$buildingSheet.Range("B2").Value = "臺北市大安區龍泉段一小段05819000建號"
$buildingSheet.Range("H2").Value = "25000000(土地建物與車位合併價）"

# ---------------------------------------------------------------------
# Sheet 4 ("債務" / debt) -- strip stray spaces/commas from text values
# ---------------------------------------------------------------------
$debtSheet = $wb.Worksheets.Item(4)
$debtSheet.Range("D2").Value = "陳麗卿新北市泰山區明志路"
$debtSheet.Range("D3").Value = "石旭松新北市泰山區明志路"
# Leading apostrophe keeps these as literal text ("6000000"/"4000000")
# instead of being auto-converted to numbers.
$debtSheet.Range("E2").Value = "'6000000"
$debtSheet.Range("E3").Value = "'4000000"
$debtSheet.Range("F2").Value = "96年02月06日"
$debtSheet.Range("F3").Value = "96年02月06日"
